# Update the LR-pairs TPM-derived metrics (Gdf1-Bmpr1a) with recomputed
# values from the new TPM script run. Only numeric metric cells change;
# the identifying columns (A-L, K, etc.) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs): receptor expression values + downstream specificity/
# edge-weight figures recomputed against the new TPM-derived receptor values.
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 0.1597850370197778
$ws.Range("R2").Value = 1.438065333178
$ws.Range("S2").Value = 0.03946985606179569
$ws.Range("T2").Value = 0.03946985606179568

# Row 3 (FAPs -> FAPs): receptor specificity shifts because the ECs receptor
# total (used in the specificity denominator) changed.
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("S3").Value = 0.5740512192416045
$ws.Range("T3").Value = 0.5740512192416045

# Row 4 (FAPs -> MuSCs): same denominator shift as row 3.
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("S4").Value = 0.3798028779254582
$ws.Range("T4").Value = 0.3798028779254582

# Row 5 (MuSCs -> ECs): ligand specificity rounding + receptor values/edge
# weights, same as row 2's ECs target.
$ws.Range("J5").Value = 0.006676046771141624
$ws.Range("M5").Value = 2.341355666666667
$ws.Range("N5").Value = 7.024067000000001
$ws.Range("O5").Value = 0.03973512964576821
$ws.Range("P5").Value = 0.0397351296457682
$ws.Range("Q5").Value = 0.001073901799111111
$ws.Range("R5").Value = 0.009665116192000001
$ws.Range("S5").Value = 0.0002652735839725247
$ws.Range("T5").Value = 0.0002652735839725246

# Row 6 (MuSCs -> FAPs): ligand specificity rounding + receptor specificity
# denominator shift.
$ws.Range("J6").Value = 0.006676046771141624
$ws.Range("O6").Value = 0.5779093692199981
$ws.Range("P6").Value = 0.5779093692199981
$ws.Range("Q6").Value = 0.01561887218844444
$ws.Range("S6").Value = 0.003858149978393661
$ws.Range("T6").Value = 0.003858149978393661

# Row 7 (MuSCs -> MuSCs): ligand specificity rounding + receptor specificity
# denominator shift.
$ws.Range("J7").Value = 0.006676046771141624
$ws.Range("O7").Value = 0.3823555011342337
$ws.Range("P7").Value = 0.3823555011342337
$ws.Range("S7").Value = 0.002552623208775438
$ws.Range("T7").Value = 0.002552623208775438
